# Apply updated cryptocurrency price/volume data to Sheet1 (generated by GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.110.16'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '1.667.13'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.84'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5103'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.75%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2636'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06417'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.54'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07419'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.15%  '
$ws.Range("D12").Value = '1.673.76'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.513'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5804'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008567'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.17'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.60%  '
$ws.Range("D17").Value = '26.164.31'
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.928'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.31%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.81'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '189.55'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.206'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '145.03'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.616'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1196'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.59'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06357'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +13.89%  '
$ws.Range("E29").Value = '  -1.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.317'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.528'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.507'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.635'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.014'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6082'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.363'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.649'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.162'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01607'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = '1.077.16'
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8614'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.010'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.15'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("D44").Value = '1.814.25'
$ws.Range("E45").Value = '  +9.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.15'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.061'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4293'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05197'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.924'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.57%  '
